# Update the "想去人数" (F column) values on the "展览" and "全部类型" sheets.
# Rows 2-10 currently hold 0; set them to the new counts from the upstream
# generator run (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$values = @(373, 106, 1578, 10, 23, 404, 138, 62, 451)

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($i = 0; $i -lt $values.Length; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 6).Value = $values[$i]
    }
}
